$d = $word.ActiveDocument

# The last paragraph in the document ends with the run "... UI часть" followed
# by the _GoBack bookmark. We split it into two paragraphs: the existing one
# (unchanged) and a new list item ("4.5 part") carrying the same list
# formatting (ListParagraph style / numId 2), inserted right after it.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Сделать свой маппер (отдельный модуль) с регистрацией сущностей."
